# Extend the "parallel" table from columns A:O to A:Q (adding P and Q),
# and fix the I/K and M/O column values for the data rows (2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14, Q1=15 with the same style as O1 ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Fix values in columns I, K, M, O for rows 2-25, and add P, Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: was 1
    $ws.Cells.Item($r, 11).Value = 1   # K: was 2
    $ws.Cells.Item($r, 13).Value = 2   # M: was 1
    $ws.Cells.Item($r, 15).Value = 1   # O: was 2
    $ws.Cells.Item($r, 16).Value = 2   # P: new
    $ws.Cells.Item($r, 17).Value = 2   # Q: new
}
